$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.017±0.002"
$ws.Range("C2").Value = "0.206±0.009"

$ws.Range("B3").Value = "0.080±0.018"
$ws.Range("C3").Value = "0.228±0.033"

$ws.Range("B4").Value = "0.701±0.187"
$ws.Range("C4").Value = "0.309±0.188"

$ws.Range("B5").Value = "0.976±0.019"
$ws.Range("C5").Value = "0.423±0.036"

$ws.Range("B6").Value = "0.882±0.104"
$ws.Range("C6").Value = "0.666±0.141"

$ws.Range("B7").Value = "0.755±0.208"
$ws.Range("C7").Value = "0.204±0.138"

$ws.Range("B8").Value = "0.008±0.002"
$ws.Range("C8").Value = "0.212±0.034"

$ws.Range("B9").Value = "0.119±0.047"
$ws.Range("C9").Value = "0.199±0.038"

$ws.Range("B10").Value = "0.722±0.072"
$ws.Range("C10").Value = "0.481±0.103"
